$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Landing Page")

# Row 4 ("Navbar"): task finished -> set an End Date (same date style as
# the Start Date column).
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 45003

# Row 5 ("Body"): grew a bit taller (wrapped description), content unchanged.
$ws.Rows(5).RowHeight = 43.2

# Row 6 ("Navbar list"): fix wording, add custom height.
$ws.Range("B6").Value = "Write minimal code to show content in Navbar"
$ws.Rows(6).RowHeight = 30

# Row 7 (new): Validation task.
$ws.Range("A7").Value = "Validation"
$ws.Range("B7").Value = "Add input field validation in login and signup form"
$ws.Range("C7").Value = "Completed"
$ws.Range("D7").Value = 45003
$ws.Range("D7").NumberFormat = "mm-dd-yy"
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 45004
$ws.Rows(7).RowHeight = 29.4

# Row 2 ("First glance"): shorten the description and mark it Completed.
$ws.Range("B2").Value = "Add a background image , food related."
$ws.Range("C2").Value = "Completed"

# Extend the conditional formatting range that used to stop at C6 so it
# also covers the newly added row 7.
$ws.Range("A6:C6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A6:C7"))

# Put the selection where the author left off.
$ws.Range("D9").Select()
